$p = $ppt.ActivePresentation

# The deck originally had two slides: a blank title slide (id 256) and the
# "Planetary Docs" class-diagram slide (id 257). The edit removes the blank
# title slide, leaving the class-diagram slide as the sole slide, and gives
# it a real title.
$p.Slides.Item(1).Delete()

$s = $p.Slides.Item(1)
$title = $s.Shapes.Item(1)
$tr = $title.TextFrame.TextRange
$tr.Text = "Class Diagram ("
$tr.InsertAfter("arrow means one-to-many)") | Out-Null
